# Append three new daily-log rows (106, 107, 108) to each of the four
# worksheets (DE_LFT_#1, DE_LFT_#2, DE_PLT_#1, DE_PLT_#2), continuing the
# existing per-day pattern found in rows 2-105.

$wb = $excel.ActiveWorkbook

# New timestamps (Excel serial date/time values) for the three appended rows.
$dates = @([double]"45892.43943287037", [double]"45893.43509259259", [double]"45894.438125")

# Per-sheet static payload (columns B-I) that repeats for every new row on
# that sheet; only column A (the timestamp) changes per row.
$sheetData = @{
    "DE_LFT_#1" = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x1C"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 284
        I = 14
    }
    "DE_LFT_#2" = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x20"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 288
        I = 14
    }
    "DE_PLT_#1" = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x73"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 115
        I = 7
    }
    "DE_PLT_#2" = @{
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x71"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 113
        I = 3
    }
}

$sheetNames = @("DE_LFT_#1", "DE_LFT_#2", "DE_PLT_#1", "DE_PLT_#2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $payload = $sheetData[$name]

    $startRow = 106
    for ($i = 0; $i -lt $dates.Length; $i++) {
        $row = $startRow + $i

        $cellA = $ws.Cells.Item($row, 1)
        $cellA.Value = $dates[$i]
        $cellA.NumberFormat = $ws.Cells.Item(105, 1).NumberFormat

        $ws.Cells.Item($row, 2).Value = $payload.B
        $ws.Cells.Item($row, 3).Value = $payload.C
        $ws.Cells.Item($row, 4).Value = $payload.D
        $ws.Cells.Item($row, 5).Value = $payload.E
        $ws.Cells.Item($row, 6).Value = $payload.F
        $ws.Cells.Item($row, 7).Value = $payload.G
        $ws.Cells.Item($row, 8).Value = $payload.H
        $ws.Cells.Item($row, 9).Value = $payload.I
    }
}
